$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '91.313.46'
$ws.Range("E2").Value = '  +1.81%  '

# Row 3
$ws.Range("D3").Value = '3.153.84'
$ws.Range("E3").Value = '  +2.58%  '

# Row 4
$ws.Range("E4").Value = '  +0.34%  '

# Row 5
$ws.Range("D5").Value = '''237.98'
$ws.Range("E5").Value = '  +0.80%  '

# Row 6
$ws.Range("D6").Value = '''617.29'
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("D7").Value = '''1.12'
$ws.Range("E7").Value = '  +6.01%  '

# Row 8
$ws.Range("D8").Value = '''0.373'
$ws.Range("E8").Value = '  +2.96%  '

# Row 9
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  -0.11%  '

# Row 10
$ws.Range("B10").Value = 'LidoStakedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D10").Value = '3.151.97'
$ws.Range("E10").Value = '  +2.49%  '

# Row 11
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").Value = '''0.738'
$ws.Range("E11").Value = '  +3.11%  '

# Row 12
$ws.Range("E12").Value = '  +2.54%  '

# Row 13
$ws.Range("D13").Value = '''0.0000245'
$ws.Range("E13").Value = '  -2.24%  '

# Row 14
$ws.Range("D14").Value = '''35.11'
$ws.Range("E14").Value = '  -0.58%  '

# Row 15
$ws.Range("D15").Value = '''5.53'
$ws.Range("E15").Value = '  +3.03%  '

# Row 16
$ws.Range("D16").Value = '91.139.65'
$ws.Range("E16").Value = '  +1.60%  '

# Row 17
$ws.Range("D17").Value = '3.736.36'
$ws.Range("E17").Value = '  +1.81%  '

# Row 18
$ws.Range("D18").Value = '3.184.66'
$ws.Range("E18").Value = '  +2.98%  '

# Row 19
$ws.Range("E19").Value = '  -2.61%  '

# Row 20
$ws.Range("D20").Value = '''15.10'
$ws.Range("E20").Value = '  +9.44%  '

# Row 21
$ws.Range("D21").Value = '''5.84'
$ws.Range("E21").Value = '  +7.96%  '

# Row 22
$ws.Range("D22").Value = '''0.0000201'
$ws.Range("E22").Value = '  -5.01%  '

# Row 23
$ws.Range("D23").Value = '''441.36'
$ws.Range("E23").Value = '  +1.82%  '

# Row 24
$ws.Range("D24").Value = '''9.11'
$ws.Range("E24").Value = '  +3.75%  '

# Row 25
$ws.Range("D25").Value = '''5.95'
$ws.Range("E25").Value = '  +6.49%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '''87.97'
$ws.Range("E26").Value = '  +1.54%  '

# Row 27
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").Value = '''11.87'
$ws.Range("E27").Value = '  +0.92%  '

# Row 29
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.07%  '

# Row 30
$ws.Range("D30").Value = '''0.232'
$ws.Range("E30").Value = '  +18.91%  '

# Row 31
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").Value = '''0.171'
$ws.Range("E31").Value = '  +9.37%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.122'
$ws.Range("E32").Value = '  +37.64%  '

# Row 33
$ws.Range("D33").Value = '''9.30'
$ws.Range("E33").Value = '  +2.61%  '

# Row 34
$ws.Range("E34").Value = '  +11.73%  '

# Row 35
$ws.Range("D35").Value = '''0.938'
$ws.Range("E35").Value = '  -6.16%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '''26.28'
$ws.Range("E36").Value = '  +2.76%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").Value = '''7.57'
$ws.Range("E37").Value = '  +6.39%  '

# Row 38
$ws.Range("D38").Value = '''505.08'
$ws.Range("E38").Value = '  +1.89%  '

# Row 39
$ws.Range("D39").Value = '''1.34'
$ws.Range("E39").Value = '  +6.76%  '

# Row 40
$ws.Range("E40").Value = '  +2.10%  '

# Row 41
$ws.Range("D41").Value = '''0.446'
$ws.Range("E41").Value = '  +12.60%  '

# Row 42
$ws.Range("D42").Value = '''3.81'
$ws.Range("E42").Value = '  +6.09%  '

# Row 43
$ws.Range("D43").Value = '''3.42'
$ws.Range("E43").Value = '  -8.03%  '

# Row 44
$ws.Range("E44").Value = '  +0.37%  '

# Row 45
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("D46").Value = '''159.68'
$ws.Range("E46").Value = '  +5.34%  '

# Row 47
$ws.Range("D47").Value = '''0.705'
$ws.Range("E47").Value = '  +4.49%  '

# Row 48
$ws.Range("E48").Value = '  +3.30%  '

# Row 49
$ws.Range("E49").Value = '  +4.18%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '''4.43'
$ws.Range("E50").Value = '  +2.15%  '

# Row 51
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '''44.09'
$ws.Range("E51").Value = '  -0.57%  '

